$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fixed (unchanged across the whole block) column values
$fixedA = 6
$fixedB = 'Mercado Mayorista Lo Valledor de Santiago'
$fixedC = 'Metropolitana'
$fixedE = 13
$fixedF = 100112043
$fixedG = 'Pepino dulce'
$fixedH = 'Cultivar IV Región'
$fixedN = '$/bandeja 18 kilos'
$fixedO = 'Provincia de Limarí'
$fixedQ = 18
$fixedR = 'Hortaliza'

# New values for columns D, I, J, K, L, M, P for rows 150..218 (69 rows total)
$newD = @(44642,44642,44642,44642,44357,44357,44357,44357,44641,44641,44641,44641,44329,44329,44329,44329,44294,44294,44294,44294,44264,44264,44264,44396,44396,44396,44396,44279,44279,44279,44279,44301,44301,44301,44385,44385,44385,44385,44236,44236,44236,44229,44229,44299,44299,44299,44312,44312,44312,44312,44399,44399,44399,44258,44258,44258,44258,44349,44349,44349,44349,44285,44285,44285,44285,44335,44335,44335,44335)
$newI = @("Especial","Primera","Segunda","Tercera","Especial","Primera","Segunda","Tercera","Especial","Primera","Segunda","Tercera","Especial","Primera","Segunda","Tercera","Especial","Primera","Segunda","Tercera","Especial","Primera","Segunda","Especial","Primera","Segunda","Tercera","Especial","Primera","Segunda","Tercera","Especial","Primera","Segunda","Especial","Primera","Segunda","Tercera","Especial","Primera","Segunda","Especial","Primera","Especial","Primera","Segunda","Especial","Primera","Segunda","Tercera","Especial","Primera","Segunda","Especial","Primera","Segunda","Tercera","Especial","Primera","Segunda","Tercera","Especial","Primera","Segunda","Tercera","Especial","Primera","Segunda","Tercera")
$newJ = @(230,370,280,150,220,330,240,130,150,260,130,90,130,290,240,80,180,320,260,55,80,150,130,140,280,240,75,250,290,180,90,980,720,600,130,180,120,60,130,150,60,120,230,640,490,290,240,280,180,70,80,240,160,70,150,125,60,120,230,180,60,480,530,220,80,150,290,220,70)
$newK = @(14000,12000,9000,6000,14000,12000,9000,7000,14000,13000,10000,8000,14000,12000,9000,6000,12000,9000,7000,5000,14000,12000,8000,15000,13000,10000,8000,11000,9000,7000,5000,11000,9000,7000,15000,13000,10000,8000,14000,12000,10000,15000,14000,10000,8000,6000,12000,9000,7000,6000,17000,15000,12000,14000,12000,10000,7000,14000,12000,10000,7000,12000,10000,8000,6000,14000,12000,9000,6000)
$newL = @(14000,12000,9000,6000,14000,12000,9000,7000,14000,13000,10000,8000,14000,12000,9000,6000,12000,9000,7000,5000,14000,12000,8000,15000,13000,10000,8000,11000,9000,7000,5000,12000,10000,8000,15000,13000,10000,8000,14000,12000,10000,15000,14000,11000,9000,7000,12000,9000,7000,6000,17000,15000,12000,14000,12000,10000,7000,14000,12000,10000,7000,13000,11000,8000,6000,14000,12000,9000,6000)
$newM = @(14000,12000,9000,6000,14000,12000,9000,7000,14000,13000,10000,8000,14000,12000,9000,6000,12000,9000,7000,5000,14000,12000,8000,15000,13000,10000,8000,11000,9000,7000,5000,11765,9764,7833,15000,13000,10000,8000,14000,12000,10000,15000,14000,10500,8469,6586,12000,9000,7000,6000,17000,15000,12000,14000,12000,10000,7000,14000,12000,10000,7000,12479,10340,8000,6000,14000,12000,9000,6000)
$newP = @(778,667,500,333,778,667,500,389,778,722,556,444,778,667,500,333,667,500,389,278,778,667,444,833,722,556,444,611,500,389,278,654,542,435,833,722,556,444,778,667,556,833,778,583,470,366,667,500,389,333,944,833,667,778,667,556,389,778,667,556,389,693,574,444,333,778,667,500,333)

for ($idx = 0; $idx -lt $newD.Length; $idx++) {
    $r = 150 + $idx

    $ws.Cells.Item($r, 1).Value2 = $fixedA
    $ws.Cells.Item($r, 2).Value2 = $fixedB
    $ws.Cells.Item($r, 3).Value2 = $fixedC
    $ws.Cells.Item($r, 4).Value2 = $newD[$idx]
    $ws.Cells.Item($r, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
    $ws.Cells.Item($r, 5).Value2 = $fixedE
    $ws.Cells.Item($r, 6).Value2 = $fixedF
    $ws.Cells.Item($r, 7).Value2 = $fixedG
    $ws.Cells.Item($r, 8).Value2 = $fixedH
    $ws.Cells.Item($r, 9).Value2 = $newI[$idx]
    $ws.Cells.Item($r, 10).Value2 = $newJ[$idx]
    $ws.Cells.Item($r, 11).Value2 = $newK[$idx]
    $ws.Cells.Item($r, 12).Value2 = $newL[$idx]
    $ws.Cells.Item($r, 13).Value2 = $newM[$idx]
    $ws.Cells.Item($r, 14).Value2 = $fixedN
    $ws.Cells.Item($r, 15).Value2 = $fixedO
    $ws.Cells.Item($r, 16).Value2 = $newP[$idx]
    $ws.Cells.Item($r, 17).Value2 = $fixedQ
    $ws.Cells.Item($r, 18).Value2 = $fixedR
}
